$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2 (shifts existing data rows 2-21 down to 3-22)
$ws.Rows.Item(2).Insert()
# The insert copies formatting from the row above (the bold header); strip
# it back off so the new data row matches the plain (unstyled) data rows.
$ws.Range("A2:C2").ClearFormats()

# Fill the newly inserted row 2 with its data
$ws.Cells.Item(2, 1).Value = -0.01545035839080799
$ws.Cells.Item(2, 2).Value = 0.2315296996384856
$ws.Cells.Item(2, 3).Value = -0.5798605158925052

# Append 9 new rows of data at the bottom (rows 23 through 31)
$newRows = @(
    @(-1.618811368942268, -0.6589505374431646, 0.2403407692909387),
    @(-0.08240008354186718, 0.4783504903316493, -3.809414207935333),
    @(-3.95973014831543, 0.9762580394744872, -4.069071769714356),
    @(-3.901577949523926, 1.771272063255311, -2.484678864479062),
    @(-1.957046031951897, -0.6577051877975557, -7.9572014808655),
    @(0.0412573218345611, -3.739429324865336, -3.584903955459609),
    @(-2.076164960861222, -2.966795355081547, -0.5832877159118695),
    @(-5.426012933254244, -0.3007338047027568, 0.185311913490301),
    @(-6.20224690437317, 0.8901370018720657, 1.72858691215514)
)

$row = 23
foreach ($r in $newRows) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $row++
}
